# Applies the "Updated symbol list" coin-price refresh (Wed Jan 25 13:51:00 UTC 2023).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'" + "302.04"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'" + "-3.97%"
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'" + "35.50"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'" + "1.03%"
$ws.Range("E3").Style = "Normal"
$ws.Range("E4").Value = "'" + "-1.28%"
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'" + "0.08021"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'" + "-1.61%"
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'" + "1.934"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'" + "-8.95%"
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value = "'" + "7.810"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'" + "-1.83%"
$ws.Range("E7").Style = "Normal"
$ws.Range("B8").Value = "BTSEToken"
$ws.Range("C8").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
$ws.Range("D8").Value = "'" + "2.985"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'" + "3.05%"
$ws.Range("E8").Style = "Normal"
$ws.Range("B9").Value = "MXToken"
$ws.Range("C9").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D9").Value = "'" + "0.9252"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'" + "-0.57%"
$ws.Range("E9").Style = "Normal"
$ws.Range("B10").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C10").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D10").Value = "'" + "0.1312"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'" + "27.92%"
$ws.Range("E10").Style = "Normal"
$ws.Range("B11").Value = "WazirX"
$ws.Range("C11").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D11").Value = "'" + "0.1856"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'" + "-4.40%"
$ws.Range("E11").Style = "Normal"
$ws.Range("B12").Value = "MandalaExchangeToken"
$ws.Range("C12").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D12").Value = "'" + "0.09274"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'" + "3.11%"
$ws.Range("E12").Style = "Normal"
$ws.Range("B13").Value = "BitrueCoin"
$ws.Range("C13").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D13").Value = "'" + "0.03394"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'" + "-8.99%"
$ws.Range("E13").Style = "Normal"
$ws.Range("B14").Value = "BitMartToken"
$ws.Range("C14").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D14").Value = "'" + "0.09873"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'" + "-0.29%"
$ws.Range("E14").Style = "Normal"
$ws.Range("B15").Value = "BitForexToken"
$ws.Range("C15").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D15").Value = "'" + "0.001389"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'" + "-3.21%"
$ws.Range("E15").Style = "Normal"
$ws.Range("B16").Value = "TigerCash"
$ws.Range("C16").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D16").Value = "'" + "0.005798"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'" + "-0.56%"
$ws.Range("E16").Style = "Normal"
$ws.Range("B17").Value = "LEO"
$ws.Range("C17").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D17").Value = "'" + "3.511"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'" + "1.20%"
$ws.Range("E17").Style = "Normal"
$ws.Range("B18").Value = "GateToken"
$ws.Range("C18").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("D18").Value = "'" + "4.062"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'" + "-2.16%"
$ws.Range("E18").Style = "Normal"
$ws.Range("E19").Value = "'" + "-0.19%"
$ws.Range("E19").Style = "Normal"
$ws.Range("E20").Value = "'" + "-2.34%"
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'" + "5.036"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'" + "-1.27%"
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'" + "0.2402"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'" + "8.31%"
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'" + "0.04489"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'" + "-1.45%"
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = "'" + "0.001213"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'" + "-2.90%"
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'" + "0.004804"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'" + "2.37%"
$ws.Range("E25").Style = "Normal"
$ws.Range("E26").Value = "'" + "-0.12%"
$ws.Range("E26").Style = "Normal"
$ws.Range("D27").Value = "'" + "0.0003003"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'" + "-33.45%"
$ws.Range("E27").Style = "Normal"
$ws.Range("D39").Value = "'" + "0.01913"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'" + "-1.55%"
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = "'" + "0.04734"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'" + "-3.14%"
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'" + "0.007352"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'" + "-3.20%"
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = "'" + "0.009645"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'" + "22.76%"
$ws.Range("E42").Style = "Normal"
$ws.Range("E43").Value = "'" + "-3.61%"
$ws.Range("E43").Style = "Normal"
$ws.Range("E44").Value = "'" + "0.57%"
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'" + "0.01088"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'" + "-7.64%"
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = "'" + "0.00006316"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'" + "-6.63%"
$ws.Range("E46").Style = "Normal"
$ws.Range("E47").Value = "'" + "-0.19%"
$ws.Range("E47").Style = "Normal"
$ws.Range("E48").Value = "'" + "-67.12%"
$ws.Range("E48").Style = "Normal"
$ws.Range("E49").Value = "'" + "-12.62%"
$ws.Range("E49").Style = "Normal"
$ws.Range("E50").Value = "'" + "-0.19%"
$ws.Range("E50").Style = "Normal"
$ws.Range("E51").Value = "'" + "-0.19%"
$ws.Range("E51").Style = "Normal"
